# Move the 5th slide (the picture-only "Figure thumbnail gr5" slide) to become
# the first slide of the deck. This shifts the original slides 1-4 down by one
# position (they become slides 2-5), matching the reordered <p:sldIdLst> in the
# target presentation.xml.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$s.MoveTo(1)
